{"js": "// The template contains two (otherwise identical) paragraphs that render\n// \"{{ impound_list }}.\" \u2014 one in the motion body, one in the order body.\n// The diff only touches the FIRST one (the motion paragraph, the one that\n// immediately follows the \"Now Comes ... to request that the Court\n// impound/keep private:\" paragraph): its Jinja2 expression\n// \"impound_list\" becomes \"impound_list.true_values()\".\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet targetParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"impound_list\") !== -1) {\n    targetParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!targetParagraph) {\n  throw new Error(\"Could not find a paragraph containing 'impound_list'.\");\n}\n\nconst matches = targetParagraph.search(\"impound_list\", { matchCase: true });\nmatches.load(\"items\");\nawait context.sync();\n\nif (matches.items.length === 0) {\n  throw new Error(\"Could not find 'impound_list' text inside the target paragraph.\");\n}\n\n// Replace \"impound_list\" with \"impound_list.true_values()\" so the merge\n// field reads \"{{ impound_list.true_values() }}.\" afterwards.\nmatches.items[0].insertText(\"impound_list.true_values()\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The template contains two (otherwise identical) paragraphs that render\n# \"{{ impound_list }}.\" -- one in the motion body, one in the order body.\n# The diff only touches the FIRST one (the motion paragraph, the one that\n# immediately follows the \"Now Comes ... to request that the Court\n# impound/keep private:\" paragraph): its Jinja2 expression\n# \"impound_list\" becomes \"impound_list.true_values()\".\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    $rng = $p.Range\n    if ($rng.Text -like \"*impound_list*\") {\n        # Scope Find/Replace to just this paragraph so only the first\n        # occurrence in the document is touched.\n        $find = $rng.Find\n        $find.ClearFormatting()\n        $find.Text = \"impound_list\"\n        $find.Replacement.Text = \"impound_list.true_values()\"\n        $find.Forward = $true\n        $find.Wrap = 0\n        $find.Format = $false\n        $find.MatchCase = $true\n        $find.MatchWholeWord = $false\n        $find.MatchWildcards = $false\n        $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n        break\n    }\n}\n"}
